$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 402
$ws.Range("I12").Value = 206.66667
$ws.Range("K12").Value = 206.66667
$ws.Range("M12").Value = -36.66667000000001
$ws.Range("H33").Value = 14494.782
$ws.Range("I33").Value = 15589.429
$ws.Range("K33").Value = 15589.429
$ws.Range("M33").Value = -15360.429
$ws.Range("H55").Value = 381.5
$ws.Range("I55").Value = 62.666668
$ws.Range("J55").Value = 700.3333
$ws.Range("K55").Value = 62.666668
$ws.Range("L55").Value = 700.3333
$ws.Range("M55").Value = 151.333332
$ws.Range("N55").Value = -1128.3333
$ws.Range("H62").Value = 2195.3333
$ws.Range("I62").Value = 2043.25
$ws.Range("K62").Value = 2043.25
$ws.Range("M62").Value = -1419.25
$ws.Range("H65").Value = 2195.3333
$ws.Range("I65").Value = 2043.25
$ws.Range("K65").Value = 10216.25
$ws.Range("M65").Value = -7096.25
$ws.Range("H106").Value = 2157.6
$ws.Range("I106").Value = 2129.3333
$ws.Range("K106").Value = 2129.3333
$ws.Range("M106").Value = -1498.3333
$ws.Range("H129").Value = 2193835.5
$ws.Range("I129").Value = 717.73334
$ws.Range("K129").Value = 2153.20002
$ws.Range("M129").Value = 2846.79998
$ws.Range("H131").Value = 6040.2856
$ws.Range("I131").Value = 6047
$ws.Range("K131").Value = 18141
$ws.Range("M131").Value = -13101
$ws.Range("H132").Value = 2221.077
$ws.Range("I132").Value = 1787.4
$ws.Range("J132").Value = 3666.6667
$ws.Range("K132").Value = 5362.200000000001
$ws.Range("L132").Value = 11000.0001
$ws.Range("M132").Value = -2832.200000000001
$ws.Range("N132").Value = -16060.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4866.7896
$ws.Range("J45").Value = 6327.1665
$ws.Range("L45").Value = 6327.1665
$ws.Range("N45").Value = -7081.1665
$ws.Range("H74").Value = 3474.0625
$ws.Range("I74").Value = 2790
$ws.Range("K74").Value = 2790
$ws.Range("M74").Value = -1916
$ws.Range("H77").Value = 3474.0625
$ws.Range("I77").Value = 2790
$ws.Range("K77").Value = 13950
$ws.Range("M77").Value = -9582
$ws.Range("H97").Value = 789.3461
$ws.Range("I97").Value = 635.5454999999999
$ws.Range("J97").Value = 1635.25
$ws.Range("K97").Value = 635.5454999999999
$ws.Range("L97").Value = 1635.25
$ws.Range("M97").Value = -139.5454999999999
$ws.Range("N97").Value = -2627.25
$ws.Range("H102").Value = 5264.3887
$ws.Range("I102").Value = 3917.2666
$ws.Range("K102").Value = 3917.2666
$ws.Range("M102").Value = -2295.2666
$ws.Range("H110").Value = 8259.421
$ws.Range("I110").Value = 6640.727
$ws.Range("K110").Value = 6640.727
$ws.Range("M110").Value = -4595.727
$ws.Range("H132").Value = 4239.4873
$ws.Range("I132").Value = 3653.7354
$ws.Range("K132").Value = 10961.2062
$ws.Range("M132").Value = -8431.206200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2770.7
$ws.Range("I86").Value = 2744.3076
$ws.Range("K86").Value = 2744.3076
$ws.Range("M86").Value = -1621.3076
$ws.Range("H89").Value = 2770.7
$ws.Range("I89").Value = 2744.3076
$ws.Range("K89").Value = 13721.538
$ws.Range("M89").Value = -8105.538
$ws.Range("H107").Value = 1700.7778
$ws.Range("I107").Value = 1700.7778
$ws.Range("K107").Value = 1700.7778
$ws.Range("M107").Value = 219.2221999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6902.0454
$ws.Range("I31").Value = 7249.6665
$ws.Range("K31").Value = 7249.6665
$ws.Range("M31").Value = -6954.6665
$ws.Range("H34").Value = 6902.0454
$ws.Range("I34").Value = 7249.6665
$ws.Range("K34").Value = 7249.6665
$ws.Range("M34").Value = -7047.6665
$ws.Range("H58").Value = 5141.6206
$ws.Range("I58").Value = 2576.5715
$ws.Range("J58").Value = 11874.875
$ws.Range("K58").Value = 2576.5715
$ws.Range("L58").Value = 11874.875
$ws.Range("M58").Value = -2373.5715
$ws.Range("N58").Value = -12280.875
$ws.Range("H99").Value = 5498.5
$ws.Range("I99").Value = 5712.5713
$ws.Range("K99").Value = 5712.5713
$ws.Range("M99").Value = -4214.5713
$ws.Range("H126").Value = 5498.5
$ws.Range("I126").Value = 5712.5713
$ws.Range("K126").Value = 17137.7139
$ws.Range("M126").Value = -14667.7139
$ws.Range("H134").Value = 6031.2915
$ws.Range("I134").Value = 4201.3335
$ws.Range("K134").Value = 12604.0005
$ws.Range("M134").Value = -10069.0005
$ws.Range("H136").Value = 5141.6206
$ws.Range("I136").Value = 2576.5715
$ws.Range("J136").Value = 11874.875
$ws.Range("K136").Value = 7729.7145
$ws.Range("L136").Value = 35624.625
$ws.Range("M136").Value = -5179.7145
$ws.Range("N136").Value = -40724.625
$ws.Range("H141").Value = 33094.43
$ws.Range("I141").Value = 24999
$ws.Range("J141").Value = 39166
$ws.Range("K141").Value = 24999
$ws.Range("L141").Value = 39166
$ws.Range("M141").Value = -19819
$ws.Range("N141").Value = -49526

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 33333934
$ws.Range("I11").Value = 100000000
$ws.Range("J11").Value = 901
$ws.Range("K11").Value = 300000000
$ws.Range("L11").Value = 2703
$ws.Range("M11").Value = -299999860
$ws.Range("N11").Value = -2983
$ws.Range("H12").Value = 1927.3334
$ws.Range("J12").Value = 2062
$ws.Range("L12").Value = 6186
$ws.Range("N12").Value = -6532
$ws.Range("H134").Value = 1299.8
$ws.Range("I134").Value = 1299.8
$ws.Range("K134").Value = 3899.4
$ws.Range("M134").Value = 1170.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4661.25
$ws.Range("I80").Value = 3497.75
$ws.Range("K80").Value = 3497.75
$ws.Range("M80").Value = -2499.75
$ws.Range("H83").Value = 4661.25
$ws.Range("I83").Value = 3497.75
$ws.Range("K83").Value = 17488.75
$ws.Range("M83").Value = -12496.75
$ws.Range("H97").Value = 4353.933
$ws.Range("I97").Value = 846.4545000000001
$ws.Range("J97").Value = 13999.5
$ws.Range("K97").Value = 846.4545000000001
$ws.Range("L97").Value = 13999.5
$ws.Range("M97").Value = -350.4545000000001
$ws.Range("N97").Value = -14991.5
$ws.Range("H102").Value = 2189.818
$ws.Range("I102").Value = 1683.8
$ws.Range("K102").Value = 1683.8
$ws.Range("M102").Value = -61.79999999999995
$ws.Range("H122").Value = 3321.4666
$ws.Range("I122").Value = 2572.6667
$ws.Range("J122").Value = 6316.6665
$ws.Range("K122").Value = 7718.000100000001
$ws.Range("L122").Value = 18949.9995
$ws.Range("M122").Value = -5268.000100000001
$ws.Range("N122").Value = -23849.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 49474.855
$ws.Range("I61").Value = 60740.234
$ws.Range("J61").Value = 1597
$ws.Range("K61").Value = 60740.234
$ws.Range("L61").Value = 1597
$ws.Range("M61").Value = -60538.234
$ws.Range("N61").Value = -2001
$ws.Range("H113").Value = 49474.855
$ws.Range("I113").Value = 60740.234
$ws.Range("J113").Value = 1597
$ws.Range("K113").Value = 60740.234
$ws.Range("L113").Value = 1597
$ws.Range("M113").Value = -58570.234
$ws.Range("N113").Value = -5937
$ws.Range("H122").Value = 3093.3242
$ws.Range("I122").Value = 3020.6365
$ws.Range("K122").Value = 9061.9095
$ws.Range("M122").Value = -6611.9095
$ws.Range("H123").Value = 68000
$ws.Range("J123").Value = 68000
$ws.Range("L123").Value = 68000
$ws.Range("N123").Value = -77800
$ws.Range("H132").Value = 8121.311
$ws.Range("I132").Value = 7954.15
$ws.Range("K132").Value = 23862.45
$ws.Range("M132").Value = -21332.45

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5227
$ws.Range("I81").Value = 4431.1665
$ws.Range("K81").Value = 8862.333000000001
$ws.Range("M81").Value = -7801.333000000001
$ws.Range("H84").Value = 5227
$ws.Range("I84").Value = 4431.1665
$ws.Range("K84").Value = 44311.665
$ws.Range("M84").Value = -39007.665
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H122").Value = 4035.5557
$ws.Range("I122").Value = 2045.7142
$ws.Range("J122").Value = 11000
$ws.Range("K122").Value = 6137.142599999999
$ws.Range("L122").Value = 33000
$ws.Range("M122").Value = -3687.142599999999
$ws.Range("N122").Value = -37900
$ws.Range("H132").Value = 3048
$ws.Range("I132").Value = 2710.0977
$ws.Range("K132").Value = 8130.293099999999
$ws.Range("M132").Value = -5600.293099999999
